$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new timestamp values to column A (rows 118-121)
$ws.Range("A118").Value = 51.3
$ws.Range("A119").Value = 51.7
$ws.Range("A120").Value = 51.9
$ws.Range("A121").Value = 52.3

# Update the active selection to reflect the new last cell
$ws.Range("A122").Select()
